$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/link/percentage updates (safe as plain text, non-numeric-looking) ---
$ws.Range('D2').Value = '68.745.07'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = '3.822.09'
$ws.Range('E3').Value = '  +3.88%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('E6').Value = '  -3.62%  '
$ws.Range('D7').Value = '3.813.54'
$ws.Range('E7').Value = '  +3.71%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('E13').Value = '  -2.54%  '
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').Value = '4.464.94'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').Value = '3.840.73'
$ws.Range('E16').Value = '  +4.28%  '
$ws.Range('D17').Value = '68.908.98'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('E18').Value = '  -1.86%  '
$ws.Range('E19').Value = '  -3.85%  '
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('E22').Value = '  +4.30%  '
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('E25').Value = '  -4.21%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('E26').Value = '  +7.43%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('E28').Value = '  -10.23%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +2.74%  '
$ws.Range('E31').Value = '  +1.91%  '
$ws.Range('E32').Value = '  +9.53%  '
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E40').Value = '  +5.97%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('D46').Value = '2.869.44'
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('E51').Value = '  -2.38%  '

# --- Numeric-looking Price values: force Text format first so Excel keeps them as strings ---
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '598.65'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.73'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.528'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.162'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.33'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.474'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '38.84'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000248'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.38'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.21'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '495.62'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.50'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.739'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '86.95'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.39'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000138'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.47'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.18'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.95'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.47'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '32.80'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.80'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.112'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.97'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.137'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.328'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '455.17'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '49.34'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.03'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.86'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.45'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.86'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0357'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.00'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '139.14'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '27.22'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.36'

# --- Strip the Text number-format back off those cells (copy plain format from an empty cell) ---
# F1 is outside the used range (A1:E51) and carries the default/general style, so copying its
# format (only) back onto the numeric-text cells clears the "@" format/quote-prefix style we
# added above while leaving the freshly-typed text VALUE untouched.
$ws.Range("F1").Copy()
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('D13').PasteSpecial(-4122)
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('D23').PasteSpecial(-4122)
$ws.Range('D24').PasteSpecial(-4122)
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('D26').PasteSpecial(-4122)
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('D31').PasteSpecial(-4122)
$ws.Range('D32').PasteSpecial(-4122)
$ws.Range('D33').PasteSpecial(-4122)
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('D38').PasteSpecial(-4122)
$ws.Range('D39').PasteSpecial(-4122)
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('D45').PasteSpecial(-4122)
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('D48').PasteSpecial(-4122)
$ws.Range('D49').PasteSpecial(-4122)
$ws.Range('D50').PasteSpecial(-4122)
$ws.Range('D51').PasteSpecial(-4122)

$excel.CutCopyMode = 0
